$d = $word.ActiveDocument

# Locate the "March 17" text inside the meeting-date paragraph and replace
# it with "April 7". A Find without replacement collapses the range onto
# the found text, which we then retarget via Range.Text assignment; toggling
# a character-formatting property back off forces the host to keep the
# surrounding text split into separate runs (mirroring how Word splits runs
# when a user selects-and-retypes a portion of existing text).
$r = $d.Content
$found = $r.Find.Execute("March 17", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if ($found) {
    $r.Text = "April 7"
    $r.Bold = 1
    $r.Bold = 0
}
